# Auto-generated Excel COM-interop script applying the Famfrit_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 262.625
$ws.Range("I12").Value = 74.25
$ws.Range("K12").Value = 74.25
$ws.Range("M12").Value = 95.75

$ws.Range("H15").Value = 810.6667
$ws.Range("I15").Value = 810.6667
$ws.Range("K15").Value = 2432.0001
$ws.Range("M15").Value = -2263.0001

$ws.Range("H19").Value = 3514.5789
$ws.Range("I19").Value = 1176.6666
$ws.Range("J19").Value = 5618.7
$ws.Range("K19").Value = 1176.6666
$ws.Range("L19").Value = 5618.7
$ws.Range("M19").Value = -1001.6666
$ws.Range("N19").Value = -5968.7

$ws.Range("H28").Value = 4394.1665
$ws.Range("I28").Value = 706.25
$ws.Range("K28").Value = 706.25
$ws.Range("M28").Value = -221.25

$ws.Range("H88").Value = 786
$ws.Range("I88").Value = 717.3333
$ws.Range("K88").Value = 717.3333
$ws.Range("M88").Value = -311.3333

$ws.Range("H91").Value = 786
$ws.Range("I91").Value = 717.3333
$ws.Range("K91").Value = 717.3333
$ws.Range("M91").Value = 686.6667

$ws.Range("H107").Value = 1649.9524
$ws.Range("I107").Value = 1172.1875
$ws.Range("K107").Value = 1172.1875
$ws.Range("M107").Value = 747.8125

$ws.Range("H138").Value = 7830.273
$ws.Range("J138").Value = 8769.666999999999
$ws.Range("L138").Value = 26309.001
$ws.Range("N138").Value = -36589.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17866104
$ws.Range("J32").Value = 24223
$ws.Range("L32").Value = 24223
$ws.Range("N32").Value = -24797

$ws.Range("H74").Value = 111238110
$ws.Range("I74").Value = 125142620
$ws.Range("K74").Value = 125142620
$ws.Range("M74").Value = -125141746

$ws.Range("H77").Value = 111238110
$ws.Range("I77").Value = 125142620
$ws.Range("K77").Value = 625713100
$ws.Range("M77").Value = -625708732

$ws.Range("H110").Value = 20175.05
$ws.Range("I110").Value = 21136.947
$ws.Range("J110").Value = 1899
$ws.Range("K110").Value = 21136.947
$ws.Range("L110").Value = 1899
$ws.Range("M110").Value = -19091.947
$ws.Range("N110").Value = -5989

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2294.0435
$ws.Range("I20").Value = 3393.8
$ws.Range("J20").Value = 1448.0769
$ws.Range("K20").Value = 3393.8
$ws.Range("L20").Value = 1448.0769
$ws.Range("M20").Value = -3146.8
$ws.Range("N20").Value = -1942.0769

$ws.Range("H50").Value = 42888.4
$ws.Range("J50").Value = 42888.4
$ws.Range("L50").Value = 42888.4
$ws.Range("N50").Value = -44036.4

$ws.Range("H86").Value = 22507.883
$ws.Range("I86").Value = 12924.25
$ws.Range("K86").Value = 12924.25
$ws.Range("M86").Value = -11801.25

$ws.Range("H89").Value = 22507.883
$ws.Range("I89").Value = 12924.25
$ws.Range("K89").Value = 64621.25
$ws.Range("M89").Value = -59005.25

$ws.Range("H105").Value = 34503
$ws.Range("I105").Value = 100010
$ws.Range("K105").Value = 100010
$ws.Range("M105").Value = -98263

$ws.Range("H109").Value = 65000
$ws.Range("J109").Value = 65000
$ws.Range("L109").Value = 65000
$ws.Range("N109").Value = -67774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2535874.8
$ws.Range("I4").Value = 4030399.5
$ws.Range("J4").Value = 45000
$ws.Range("K4").Value = 4030399.5
$ws.Range("L4").Value = 45000
$ws.Range("M4").Value = -4030287.5
$ws.Range("N4").Value = -45224

$ws.Range("H16").Value = 1854.3684
$ws.Range("I16").Value = 1907.4286
$ws.Range("K16").Value = 1907.4286
$ws.Range("M16").Value = -1620.4286

$ws.Range("H48").Value = 34733.75
$ws.Range("J48").Value = 34733.75
$ws.Range("L48").Value = 34733.75
$ws.Range("N48").Value = -35685.75

$ws.Range("H54").Value = 33595.4
$ws.Range("J54").Value = 32666.334
$ws.Range("L54").Value = 32666.334
$ws.Range("N54").Value = -33982.334

$ws.Range("H62").Value = 7916.5
$ws.Range("I62").Value = 7999
$ws.Range("J62").Value = 7900
$ws.Range("K62").Value = 7999
$ws.Range("L62").Value = 7900
$ws.Range("M62").Value = -7375
$ws.Range("N62").Value = -9148

$ws.Range("H65").Value = 7916.5
$ws.Range("I65").Value = 7999
$ws.Range("J65").Value = 7900
$ws.Range("K65").Value = 39995
$ws.Range("L65").Value = 39500
$ws.Range("M65").Value = -36875
$ws.Range("N65").Value = -45740

$ws.Range("H86").Value = 5742.143
$ws.Range("J86").Value = 6157.8
$ws.Range("L86").Value = 6157.8
$ws.Range("N86").Value = -8403.799999999999

$ws.Range("H89").Value = 5742.143
$ws.Range("J89").Value = 6157.8
$ws.Range("L89").Value = 30789
$ws.Range("N89").Value = -42021

$ws.Range("H113").Value = 1854.3684
$ws.Range("I113").Value = 1907.4286
$ws.Range("K113").Value = 1907.4286
$ws.Range("M113").Value = 262.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1903616.1
$ws.Range("I2").Value = 1195.7142
$ws.Range("K2").Value = 7174.285199999999
$ws.Range("M2").Value = -7061.285199999999

$ws.Range("H25").Value = 1977.625
$ws.Range("I25").Value = 785
$ws.Range("J25").Value = 3965.3333
$ws.Range("K25").Value = 2355
$ws.Range("L25").Value = 11895.9999
$ws.Range("M25").Value = -2186
$ws.Range("N25").Value = -12233.9999

$ws.Range("H30").Value = 1977.625
$ws.Range("I30").Value = 785
$ws.Range("J30").Value = 3965.3333
$ws.Range("K30").Value = 2355
$ws.Range("L30").Value = 11895.9999
$ws.Range("M30").Value = -2253
$ws.Range("N30").Value = -12099.9999

$ws.Range("H33").Value = 105.73333
$ws.Range("J33").Value = 131.25
$ws.Range("L33").Value = 787.5
$ws.Range("N33").Value = -1353.5

$ws.Range("H131").Value = 50606.434
$ws.Range("J131").Value = 8136.2856
$ws.Range("L131").Value = 24408.8568
$ws.Range("N131").Value = -34488.8568

$ws.Range("H133").Value = 15979.4
$ws.Range("J133").Value = 19956.334
$ws.Range("L133").Value = 59869.00199999999
$ws.Range("N133").Value = -69989.00199999999

$ws.Range("H140").Value = 2271.7778
$ws.Range("I140").Value = 1814.1111
$ws.Range("J140").Value = 2729.4443
$ws.Range("K140").Value = 5442.3333
$ws.Range("L140").Value = 8188.3329
$ws.Range("M140").Value = -262.3333000000002
$ws.Range("N140").Value = -18548.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 46029.5
$ws.Range("I32").Value = 44999.5
$ws.Range("J32").Value = 48089.5
$ws.Range("K32").Value = 44999.5
$ws.Range("L32").Value = 48089.5
$ws.Range("M32").Value = -44703.5
$ws.Range("N32").Value = -48681.5

$ws.Range("H44").Value = 25000
$ws.Range("J44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -26192

$ws.Range("H70").Value = 5058.552
$ws.Range("I70").Value = 4768.227
$ws.Range("K70").Value = 4768.227
$ws.Range("M70").Value = -4498.227

$ws.Range("H73").Value = 5058.552
$ws.Range("I73").Value = 4768.227
$ws.Range("K73").Value = 4768.227
$ws.Range("M73").Value = -3832.227

$ws.Range("H97").Value = 2670
$ws.Range("I97").Value = 1147.2858
$ws.Range("K97").Value = 1147.2858
$ws.Range("M97").Value = -651.2858000000001

$ws.Range("H126").Value = 4922367.5
$ws.Range("I126").Value = 2385101.8
$ws.Range("J126").Value = 10531061
$ws.Range("K126").Value = 7155305.399999999
$ws.Range("L126").Value = 31593183
$ws.Range("M126").Value = -7152835.399999999
$ws.Range("N126").Value = -31598123

$ws.Range("H132").Value = 2335.4883
$ws.Range("I132").Value = 2174.8572
$ws.Range("J132").Value = 3038.25
$ws.Range("K132").Value = 6524.571599999999
$ws.Range("L132").Value = 9114.75
$ws.Range("M132").Value = -3994.571599999999
$ws.Range("N132").Value = -14174.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 634.55554
$ws.Range("I55").Value = 331.7
$ws.Range("J55").Value = 1013.125
$ws.Range("K55").Value = 331.7
$ws.Range("L55").Value = 1013.125
$ws.Range("M55").Value = -158.7
$ws.Range("N55").Value = -1359.125

$ws.Range("H68").Value = 5999

$ws.Range("H71").Value = 5999

$ws.Range("H93").Value = 1752.4615
$ws.Range("I93").Value = 1434.7273
$ws.Range("J93").Value = 3500
$ws.Range("K93").Value = 1434.7273
$ws.Range("L93").Value = 3500
$ws.Range("M93").Value = -186.7273
$ws.Range("N93").Value = -5996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J122").Value = 2910.3572
$ws.Range("L122").Value = 8731.071599999999
$ws.Range("N122").Value = -13631.0716

$ws.Range("H126").Value = 6152.579
$ws.Range("I126").Value = 6187.4375
$ws.Range("K126").Value = 18562.3125
$ws.Range("M126").Value = -16092.3125

$ws.Range("H127").Value = 27000
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H136").Value = 1696.4762
$ws.Range("I136").Value = 1058.8667
$ws.Range("K136").Value = 3176.6001
$ws.Range("M136").Value = -626.6001000000001
